$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Multilingual_spreadsheet_sample")
$ws1.Range("E5").Value = "Переклад 4"
$ws1.Range("E6").Value = "Переклад 5"

$ws2 = $wb.Worksheets.Item("Sheet1")
$ws2.Range("E5").Value = "Переклад 4"
$ws2.Range("E6").Value = "Переклад 5"
